# "pop is entirely derived" — the Atom table's `pop` column (C12:C16,
# including its header "pop" in C12) is no longer stored data; it is
# computed elsewhere, so we remove it from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atoms")
$ws.Activate()

# Select the column of values being removed (mirrors the author's manual
# selection before deleting) and clear it out. Using ClearContents (not a
# shifting Delete) leaves the rest of the table - including column D on
# other rows - untouched.
$rng = $ws.Range("C12:C16")
$rng.Select()
$rng.ClearContents()
